$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from "Planification Initiale" to "Planning"
$ws.Name = "Planning"

# Update the Gantt bar for the "Analyse temps des SLAs" row (row 42):
# the Friday cell (W42) now shows "3h15" instead of "5h15", using the
# filled/blue "bar" style already used by the neighbouring bar cells.
$fmtSrc = $ws.Range("W47")
$fmtSrc.Copy()
$w42 = $ws.Range("W42")
$w42.PasteSpecial(-4122)
$w42.Value = "3h15"

# Fill in the previously empty Friday cell for the "Manuel de mise en
# service" row (row 46) with "0h45", again reusing the filled bar style.
$w46 = $ws.Range("W46")
$fmtSrc.Copy()
$w46.PasteSpecial(-4122)
$w46.Value = "0h45"

# Zoom the sheet view to 85% and leave the selection on Y42, matching the
# state the workbook was saved in.
$excel.ActiveWindow.Zoom = 85
$ws.Range("Y42").Select()
